$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.686.20'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '3.421.56'
$ws.Range("E3").Value = '  -1.64%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.51'
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.60'
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.604'
$ws.Range("E7").Value = '  +3.18%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '3.424.21'
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.14'
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").Value = '4.010.78'
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  -3.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.58'
$ws.Range("E16").Value = '  -3.61%  '
$ws.Range("D17").Value = '64.671.48'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").Value = '3.379.65'
$ws.Range("E18").Value = '  -3.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.35'
$ws.Range("E19").Value = '  -1.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.87'
$ws.Range("E20").Value = '  -2.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.52'
$ws.Range("E21").Value = '  -2.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.00'
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.33'
$ws.Range("E25").Value = '  -1.98%  '
$ws.Range("E26").Value = '  -4.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.07'
$ws.Range("E27").Value = '  +5.27%  '
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  +2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.20'
$ws.Range("E31").Value = '  -3.89%  '
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.20'
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.08'
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("E35").Value = '  +4.02%  '
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.93'
$ws.Range("E37").Value = '  -2.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0757'
$ws.Range("E38").Value = '  -2.17%  '
$ws.Range("D39").Value = '2.884.62'
$ws.Range("E39").Value = '  -5.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.76'
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.44'
$ws.Range("E41").Value = '  -3.17%  '
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.03'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E44").Value = '  -2.23%  '
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.77'
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '320.80'
$ws.Range("E47").Value = '  +3.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.23'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("E49").Value = '  -5.16%  '
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("E51").Value = '  -2.67%  '
